# Sprint Backlog updates: add new "methods" column (D) detail text for the
# session-memento related stories, tidy up a trailing space in an existing
# cell, drop the stray empty/no-op styled cell in J4, and refresh the
# selection/row-height bookkeeping that Excel would normally update when a
# user edits these rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint Backlog")

# Make sure the Sprint Backlog sheet stays the active tab (it was already
# the active tab in the source workbook).
$ws.Activate()

# Row 4 had a leftover styled-but-empty cell in J4 with no value - remove it.
$ws.Range("J4").Clear()

# Row 5: "Session save & load" story gains a clarified description.
$ws.Range("D5").Value = "Load and save a Session"

# Row 6: "Session memento" story gains its method/description.
$ws.Range("D6").Value = "implement memento Pattern"

# Row 7: "Patient selection" story gains its method/description.
$ws.Range("D7").Value = "Show patient information"

# Row 8: "Patient session information" story gains its method/description.
$ws.Range("D8").Value = "Show patient session information"

# Row 9: drop the trailing space in the existing description.
$ws.Range("D9").Value = "Prepare Presentation for Memento pattern"

# Column widths shifted slightly across the whole sheet (font-metric driven
# autosize), most notably column D which needed to widen a lot to fit the
# new, longer text.
$ws.Columns.Item(1).ColumnWidth = 5.405882352941177
$ws.Columns.Item(2).ColumnWidth = 6.252941176470587
$ws.Columns.Item(3).ColumnWidth = 21.82549019607847
$ws.Columns.Item(4).ColumnWidth = 38.44117647058827
$ws.Columns.Item(5).ColumnWidth = 13.170588235294167
$ws.Columns.Item(6).ColumnWidth = 11.656862745098065
$ws.Columns.Item(7).ColumnWidth = 9.801960784313765
$ws.Columns.Item(8).ColumnWidth = 6.088235294117647
$ws.Columns.Item(9).ColumnWidth = 10.927450980392166
$ws.Columns.Item(10).ColumnWidth = 12.276470588235266
$ws.Columns.Item(11).ColumnWidth = 16.21372549019607

# Re-assert the existing row heights so the workbook keeps remembering them
# as explicit/custom (matches the source workbook's row metadata).
$ws.Rows.Item(1).RowHeight = 14
$ws.Rows.Item(3).RowHeight = 14
$ws.Rows.Item(4).RowHeight = 14
$ws.Rows.Item(5).RowHeight = 14
$ws.Rows.Item(6).RowHeight = 14
$ws.Rows.Item(7).RowHeight = 14.9
$ws.Rows.Item(8).RowHeight = 14.9
$ws.Rows.Item(9).RowHeight = 13.3
$ws.Rows.Item(10).RowHeight = 13.3

# Selection cursor ends up on D9 after the edits.
$ws.Range("D9").Select()
